$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = -0.3623658873974311
$ws.Range("E2").Value = 0.1825419310453658

# Row 3
$ws.Range("C3").Value = -0.009261555895478946
$ws.Range("E3").Value = 0.1145211022186787

# Row 4
$ws.Range("C4").Value = -1.404263945418582
$ws.Range("E4").Value = -0.807808220045203

# Row 5
$ws.Range("C5").Value = 1.692932643509848
$ws.Range("E5").Value = 0.6262577107155831

# Row 6
$ws.Range("C6").Value = 1.020829760720643
$ws.Range("E6").Value = 1.148272834981245

# Row 7
$ws.Range("C7").Value = 0.6772121200332215
$ws.Range("E7").Value = 1.258913537332895

# Row 8
$ws.Range("C8").Value = 1.019715257608911
$ws.Range("E8").Value = 0.9536145745415947

# Row 9
$ws.Range("C9").Value = 2.173959184500385
$ws.Range("E9").Value = 1.566646323486043

# Row 10
$ws.Range("C10").Value = 1.707434489469994
$ws.Range("E10").Value = 1.30258347990615

# Row 11
$ws.Range("C11").Value = 1.456988786619839
$ws.Range("E11").Value = 1.842797144428188

# Row 12
$ws.Range("C12").Value = 1.241332692055597
$ws.Range("E12").Value = 1.58004210678635

# Row 13
$ws.Range("C13").Value = 1.592885137608979
$ws.Range("E13").Value = 1.604795846351514

# Row 14
$ws.Range("C14").Value = -2.015335584265165
$ws.Range("E14").Value = -1.215549235925828

# Row 15
$ws.Range("C15").Value = -3.579597300369253
$ws.Range("E15").Value = -1.403103901755631

# Row 16
$ws.Range("C16").Value = 4.461954539041502
$ws.Range("E16").Value = 0.7797949948739058

# Row 17
$ws.Range("C17").Value = -1.305206755692701
$ws.Range("E17").Value = 0.5821000732047832

# Row 18
$ws.Range("C18").Value = 0.0845726262934221
$ws.Range("E18").Value = 0.1341520870597357

# Row 19
$ws.Range("C19").Value = 0.9724700385226326
$ws.Range("E19").Value = 0.6236501628417823
